$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date-cell style (A301 uses style index 2: bordered, bold,
# centered, datetime number format) down through the new A-column cells so the
# appended rows keep the same look as the rest of the date column.
$ws.Range("A301").Copy($ws.Range("A302:A328"))

# New daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile per 100k ab.)
$data = @(
    @(44376, 0, 15, 85.37765382207297),
    @(44377, 0, 14, 79.68581023393477),
    @(44378, 1, 9, 51.22659229324378),
    @(44379, 0, 9, 51.22659229324378),
    @(44380, 0, 7, 39.84290511696739),
    @(44381, 0, 1, 5.691843588138198),
    @(44382, 0, 1, 5.691843588138198),
    @(44383, 0, 1, 5.691843588138198),
    @(44384, 0, 1, 5.691843588138198),
    @(44385, 0, 0, 0),
    @(44386, 2, 2, 11.3836871762764),
    @(44387, 0, 2, 11.3836871762764),
    @(44388, 0, 2, 11.3836871762764),
    @(44389, 3, 5, 28.45921794069099),
    @(44390, 2, 7, 39.84290511696739),
    @(44391, 0, 7, 39.84290511696739),
    @(44392, 0, 7, 39.84290511696739),
    @(44393, 0, 5, 28.45921794069099),
    @(44394, 0, 5, 28.45921794069099),
    @(44395, 0, 5, 28.45921794069099),
    @(44396, 1, 3, 17.07553076441459),
    @(44397, 0, 1, 5.691843588138198),
    @(44398, 3, 4, 22.76737435255279),
    @(44399, 3, 7, 39.84290511696739),
    @(44400, 0, 7, 39.84290511696739),
    @(44401, 2, 9, 51.22659229324378),
    @(44402, 1, 10, 56.91843588138198)
)

$row = 302
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
